# Commit before pulling app changes
# Adds KDS/Expo/Receipt/SandSidesDesserts/Protein display-name columns (L:P)
# plus SidesValidation / ProteinValidation columns (Q:R) to the NCR Kiosk
# test-case sheet, fills in row 2 & 3 data, and fixes the Sides value for
# row 3 (Fried Pickles -> Cajun Fries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header row (row 1), columns L..R ------------------------------
$ws.Cells.Item(1, 12).Value = "KDSMealDisplayName"
$ws.Cells.Item(1, 13).Value = "ExpoLabelDisplayName"
$ws.Cells.Item(1, 14).Value = "ReceiptDisplayName"
$ws.Cells.Item(1, 15).Value = "SandSidesDessertsDisplayName"
$ws.Cells.Item(1, 16).Value = "ProteinDisplayName"
$ws.Cells.Item(1, 17).Value = "SidesValidation"
$ws.Cells.Item(1, 18).Value = "ProteinValidation"

# ---- Row 2 (1/4 lb Popcorn Shrimp Combo) new columns L..R ---------------
$ws.Cells.Item(2, 12).Value = "Shrmp Pop 1/4LB Cmb"
$ws.Cells.Item(2, 13).Value = "Shr Pop 1/4Lb, Fries Reg"
$ws.Cells.Item(2, 14).Value = "1/4 lb Popcorn Shrimp Combo, 1/4 lb Popcorn Shrimp, Rg Cajun Fries"
$ws.Cells.Item(2, 15).Value = "Fries Reg"
$ws.Cells.Item(2, 16).Value = "Shr Pop 1/4Lb"
$ws.Cells.Item(2, 17).Value = "YES"
$ws.Cells.Item(2, 18).Value = "YES"

# ---- Row 3 (6 Pc Boneless Wings Combo) new columns L..R ------------------
$ws.Cells.Item(3, 12).Value = "Boneless 6P Combo"
$ws.Cells.Item(3, 13).Value = "Bnls 6P Classic, Fries Reg"
$ws.Cells.Item(3, 14).Value = "6 Pc. Boneless Wings Combo, 6P Classic Boneless, Rg Cajun Fries"
$ws.Cells.Item(3, 15).Value = "Fries Reg"
$ws.Cells.Item(3, 16).Value = "Bnls 6P Classic"
$ws.Cells.Item(3, 17).Value = "YES"
$ws.Cells.Item(3, 18).Value = "YES"

# ---- Fix Sides for row 3: Fried Pickles -> Cajun Fries -------------------
$ws.Cells.Item(3, 8).Value = "Cajun Fries"

# ---- Column widths for the new columns (L..Q) -----------------------------
# ColumnWidth is specified in characters; the saved OOXML width carries a
# fixed +0.8333... padding on top of it, so back that out here.
$ws.Columns.Item(12).ColumnWidth = 23.666666666666668
$ws.Columns.Item(13).ColumnWidth = 27.67666666666667
$ws.Columns.Item(14).ColumnWidth = 48.60666666666666
$ws.Columns.Item(15).ColumnWidth = 27.816666666666666
$ws.Columns.Item(16).ColumnWidth = 15.646666666666667
$ws.Columns.Item(17).ColumnWidth = 15.796666666666665

# ---- Selection matches the end state of the edit --------------------------
$ws.Range("M3").Select() | Out-Null
